# UPDATE Practico de handicap
#
# The document's visible content lives entirely inside legacy VML
# (w:pict / v:textbox) text boxes, which this host's Range/Find/Shapes
# object model does not expose as navigable paragraphs/runs. The
# reliable way to edit such content through the Word COM surface here
# is to round-trip the whole package through Content.WordOpenXML /
# Content.InsertXML, after making the precise textual edits described
# by the change:
#
#   1. Drop the now-stray w:proofErr spell-check markers around the
#      "Twitter" and "Facebook" runs.
#   2. Drop the w:proofErr grammar-check markers that bracketed the
#      lone "n" run inside the "(0,n)" and "(1,n)" callouts, and merge
#      the three split runs "(0," + "n" + ")" (and "(1," + "n" + ")")
#      back into a single run each.

$d = $word.ActiveDocument

$xml = $d.Content.WordOpenXML

# 1) Remove proofErr spellcheck markers (Twitter / Facebook) and any
#    other stray proofErr markers (gramStart/gramEnd included) in one
#    general pass - Word itself no longer needs these once the text is
#    corrected/merged, and the target revision drops all of them.
$xml = $xml -replace '<w:proofErr[^/]*/>', ''

# 2) Merge the split "(0,n)" run back into a single run.
$xml = $xml.Replace(
    '<w:r><w:t>(0,</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>)</w:t></w:r>',
    '<w:r><w:t>(0,n)</w:t></w:r>'
)

# 3) Merge the split "(1,n)" run back into a single run.
$xml = $xml.Replace(
    '<w:r><w:t>(1,</w:t></w:r><w:r><w:t>n</w:t></w:r><w:r><w:t>)</w:t></w:r>',
    '<w:r><w:t>(1,n)</w:t></w:r>'
)

$d.Content.InsertXML($xml)
